$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the Kyrgyz title in A1 with corrected wording.
$ws.Range("A1").Value = "3.9.2 Коопсуздук суунун, коопсуздук санитариянын жана гигиенанын жоктугунан болгон өлүм"

# 2. Add the new 2022 data column (S), copying the formatting from column R
#    (the previous last data column) for each affected row.
$ws.Range("R4:R14").Copy($ws.Range("S4:S14")) | Out-Null

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 1.2
$ws.Range("S6").Value = 2.7
$ws.Range("S7").Value = 0.9
$ws.Range("S8").Value = 0.4
$ws.Range("S9").Value = 0.7
$ws.Range("S10").Value = 0.9
$ws.Range("S11").Value = 1.1
$ws.Range("S12").Value = 2.7
$ws.Range("S13").Value = 0.4
$ws.Range("S14").Value = 0.6

# Make sure the header / top-left cell is selected (clears any stray selection
# like S17 left over from before the new column existed).
$ws.Range("A1").Select() | Out-Null

Write-Output "done"
